$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab from "Contact Success" to "Contact"
$ws.Name = "Contact"

# Add a "Created" value into column A for each data row (2-8), leaving the
# header row (A1 = "Status") untouched.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = "Created"
}

$ws.Range("D19").Select() | Out-Null
